$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> target D, J, K, L, M, P values
$rowData = @{
    2 = @{ D = 44222; J = 45; K = 3000; L = 3000; M = 3000; P = 1000 }
    3 = @{ D = 44557; J = 104; K = 2000; L = 2500; M = 2260; P = 753 }
    4 = @{ D = 44243; J = 45; K = 3000; L = 3000; M = 3000; P = 1000 }
    5 = @{ D = 44225; J = 56; K = 3000; L = 3000; M = 3000; P = 1000 }
    6 = @{ D = 44187; J = 65; K = 3000; L = 3000; M = 3000; P = 1000 }
    7 = @{ D = 44224; J = 67; K = 3000; L = 3000; M = 3000; P = 1000 }
    8 = @{ D = 44390; J = 50; K = 3000; L = 3000; M = 3000; P = 1000 }
    9 = @{ D = 44574; J = 50; K = 3000; L = 3000; M = 3000; P = 1000 }
    10 = @{ D = 44242; J = 95; K = 2500; L = 3000; M = 2737; P = 912 }
    11 = @{ D = 44389; J = 81; K = 2800; L = 3000; M = 2889; P = 963 }
    12 = @{ D = 44221; J = 50; K = 2500; L = 2500; M = 2500; P = 833 }
    13 = @{ D = 44260; J = 60; K = 3500; L = 3500; M = 3500; P = 1167 }
    14 = @{ D = 44292; J = 40; K = 3000; L = 3000; M = 3000; P = 1000 }
    16 = @{ D = 44291; J = 45; K = 3000; L = 3000; M = 3000; P = 1000 }
    17 = @{ D = 44536; J = 125; K = 2200; L = 2200; M = 2200; P = 733 }
    18 = @{ D = 44340; J = 54; K = 3000; L = 3000; M = 3000; P = 1000 }
    19 = @{ D = 44627; J = 78; K = 3500; L = 3500; M = 3500; P = 1167 }
    20 = @{ D = 44165; J = 68; K = 3000; L = 3000; M = 3000; P = 1000 }
    21 = @{ D = 44179; J = 78; K = 3000; L = 3000; M = 3000; P = 1000 }
    22 = @{ D = 44193; J = 70; K = 3000; L = 3000; M = 3000; P = 1000 }
    23 = @{ D = 44537; J = 88; K = 2000; L = 2200; M = 2091; P = 697 }
    24 = @{ D = 44559; J = 68; K = 2000; L = 2000; M = 2000; P = 667 }
    25 = @{ D = 44223; J = 80; K = 2500; L = 3000; M = 2781; P = 927 }
}

foreach ($row in $rowData.Keys) {
    $v = $rowData[$row]
    $ws.Range("D$row").Value = $v.D
    $ws.Range("J$row").Value = $v.J
    $ws.Range("K$row").Value = $v.K
    $ws.Range("L$row").Value = $v.L
    $ws.Range("M$row").Value = $v.M
    $ws.Range("P$row").Value = $v.P
}

Write-Host "Updated $($rowData.Count) rows"
